$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml): F2 418 -> 419, F3 2625 -> 2639
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 419
$wsExhibit.Range("F3").Value = 2639

# Sheet "全部类型" (sheet4.xml): F2 418 -> 419, F7 2625 -> 2639
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 419
$wsAll.Range("F7").Value = 2639
